# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2304
#   *_new  -> *_FV2310
# Also wraps the data range in an Excel Table and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Segmentname_FV2304",
  "Segmentgruppe_FV2304",
  "Segment_FV2304",
  "Datenelement_FV2304",
  "Segment ID_FV2304",
  "Code_FV2304",
  "Qualifier_FV2304",
  "Beschreibung_FV2304",
  "Bedingungsausdruck_FV2304",
  "Bedingung_FV2304",
  "diff",
  "Segmentname_FV2310",
  "Segmentgruppe_FV2310",
  "Segment_FV2310",
  "Datenelement_FV2310",
  "Segment ID_FV2310",
  "Code_FV2310",
  "Qualifier_FV2310",
  "Beschreibung_FV2310",
  "Bedingungsausdruck_FV2310",
  "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the used range into an Excel table (adds xl/tables/table1.xml +
# the <tableParts> reference on the worksheet).
$dataRange = $ws.UsedRange
$lo = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = $null

# Freeze the header row (pane split after row 1).
[void]$ws.Cells.Item(2, 1).Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "ok"
